# Actualización automática 2025-06-26 14:31:13
$wb = $excel.ActiveWorkbook

$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentaMensual    = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento    = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Sheet "VENTAS POR GRUPO": LAVABOS value for RIOS CARRION ANGEL BENIGNO (row 4)
$wsVentasPorGrupo.Range("I4").Value = 29.29

# Sheet "VENTA MENSUAL": junio (June) column for RIOS CARRION ANGEL BENIGNO (row 4)
$wsVentaMensual.Range("F4").Value = 777.8099999999999

# Sheet "VENTA MENSUAL": junio (June) TOTAL row (row 19)
$wsVentaMensual.Range("F19").Value = 22179.48

# Sheet "CUMPLIMIENTO MENSUAL": LAVABOS row (row 8) - VENTA, POR CUMPLIR, CUMPLIMIENTO
$wsCumplimiento.Range("D8").Value = 29.29
$wsCumplimiento.Range("E8").Value = 595.71
$wsCumplimiento.Range("F8").Value = 0.046864

# Sheet "CUMPLIMIENTO MENSUAL": TOTAL row (row 19) - VENTA, POR CUMPLIR, CUMPLIMIENTO
$wsCumplimiento.Range("D19").Value = 22179.48
$wsCumplimiento.Range("E19").Value = 25039.82386304603
$wsCumplimiento.Range("F19").Value = 0.4697121343493106
